# Insert a new data row at row 859 (pushing the existing rows 859-904 down
# to 860-905) and populate it with the new record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("859:859").Insert()

$ws.Range("A859").Value = 11
$ws.Range("B859").Value = "Vega Monumental Concepción"
$ws.Range("C859").Value = "Bíobío"
$ws.Range("D859").Value = 44931
$ws.Range("E859").Value = 8
$ws.Range("F859").Value = 100112033
$ws.Range("G859").Value = "Lechuga"
$ws.Range("H859").Value = "Escarola"
$ws.Range("I859").Value = "Primera"
$ws.Range("J859").Value = 250
$ws.Range("K859").Value = 9000
$ws.Range("L859").Value = 10000
$ws.Range("M859").Value = 9600
$ws.Range("N859").Value = "$/caja 15 unidades"
$ws.Range("O859").Value = "Región de Coquimbo"
$ws.Range("P859").Value = 640
$ws.Range("Q859").Value = 15
$ws.Range("R859").Value = "Hortaliza"
